# Applies the Pmax Sample 2D 5Hz normalization edit:
#  - Header: A1 "File Name" -> "Loc", B1 "Unnormalized P_max" -> "P_max"
#  - Column C ("Electrode Locations") removed entirely; column A now
#    holds just the electrode location codes (e.g. "A1") instead of the
#    original file names (e.g. "A1_bipolar_20V_5Hz.txt")
#  - A previously-missing measurement for electrode "F10" (63.496071) is
#    inserted between "F8" and "G5"
#  - Final data range becomes A1:B60 (59 data rows + header)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final (Loc, P_max) pairs in their new row order (row 2 .. row 60)
$data = @(
    @("A1", 125.587219),
    @("A7", 74.995852),
    @("A9", 67.535225),
    @("A13", 53.072186),
    @("B3", 76.025322),
    @("B6", 72.230414),
    @("B12", 56.425028),
    @("C1", 74.735457),
    @("C3", 70.429852),
    @("C7", 68.582862),
    @("C9", 65.575196),
    @("C11", 58.49406),
    @("C13", 54.15212),
    @("D6", 66.810559),
    @("D15", 53.641423),
    @("E1", 76.505741),
    @("E5", 64.862641),
    @("E11", 52.448449),
    @("E13", 58.689861),
    @("E15", 53.877595),
    @("F3", 70.395536),
    @("F8", 66.324084),
    @("F10", 63.496071),
    @("G5", 55.597415),
    @("G7", 52.668473),
    @("G9", 65.270392),
    @("G15", 55.7589),
    @("H1", 77.539248),
    @("H3", 70.135141),
    @("H14", 59.539678),
    @("I1", 50.968839),
    @("I5", 70.884029),
    @("I9", 60.613556),
    @("I13", 58.36689),
    @("I15", 56.057648),
    @("J3", 69.533607),
    @("J8", 32.341777),
    @("J11", 55.256277),
    @("K6", 68.764533),
    @("K8", 60.135155),
    @("K13", 64.01888),
    @("K15", 59.515455),
    @("L1", 76.624836),
    @("L3", 71.915517),
    @("L11", 67.682581),
    @("M1", 77.547322),
    @("M5", 72.385844),
    @("M7", 72.486772),
    @("M13", 63.510201),
    @("N4", 76.243327),
    @("N9", 71.681364),
    @("N12", 68.211445),
    @("N15", 62.060869),
    @("O1", 79.525518),
    @("O3", 83.663582),
    @("O6", 71.182777),
    @("O9", 86.168625),
    @("O11", 76.998271),
    @("O13", 70.490409)
)

# Old sheet had data through row 59 and a column C; clear that whole
# used range first so nothing stale is left behind once column C is
# dropped and an extra row is introduced.
$oldUsed = $ws.UsedRange
$oldUsed.Clear() | Out-Null

# Rewrite header
$ws.Cells.Item(1, 1).Value = "Loc"
$ws.Cells.Item(1, 2).Value = "P_max"

# Re-apply the header style (bold / border / centered) that the original
# cells had, since Clear() above wiped formatting too.
$ws.Range("A1:B1").Font.Bold = $true
$ws.Range("A1:B1").HorizontalAlignment = -4108  # xlCenter
$ws.Range("A1:B1").VerticalAlignment = -4160    # xlVAlignTop
$ws.Range("A1:B1").Borders.LineStyle = 1
$ws.Range("A1:B1").Borders.Weight = 2

# Write out the data rows starting at row 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

$wb.Save()
